# Remove heading from sequence concept slider.
#
# Each of the 7 slides in this deck has a "Program Execution and Sequence"
# heading textbox. On slides 1-5 the whole textbox shape is removed; on
# slides 6-7 only the text run is removed, leaving the now-empty paragraph
# (with its endParaRPr) in place.

$p = $ppt.ActivePresentation

$headingText = "Program Execution and Sequence"

# Slides where the whole heading textbox shape should be deleted.
$slidesToDeleteShape = @(1, 2, 3, 4, 5)
# Slides where only the run's text should be removed, keeping the shape
# and its now-empty paragraph (endParaRPr survives).
$slidesToClearRun = @(6, 7)

for ($idx = 1; $idx -le $p.Slides.Count; $idx++) {
    $s = $p.Slides.Item($idx)
    for ($shi = $s.Shapes.Count; $shi -ge 1; $shi--) {
        $sh = $s.Shapes.Item($shi)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText -and $sh.TextFrame.TextRange.Text -eq $headingText) {
            if ($slidesToDeleteShape -contains $idx) {
                $sh.Delete()
            } elseif ($slidesToClearRun -contains $idx) {
                $tr = $sh.TextFrame.TextRange
                $tr.Characters(1, $tr.Length).Delete()
            }
        }
    }
}

Write-Output "Heading removed from slides 1-7."
